# BA_tasks.xlsx edit: merge Eigenfaces/Cluster-Algo notes into B26 and remove the
# now-redundant D29/D30 TODO rows, per commit "feat: added BayesianGaussianMixture cluster comp"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update B26 ("Task" column for 2023-09-02) with the merged / extended note that now
#    also mentions the Variational Bayesian Mixture Model clustering work.
$ws.Range("B26").Value() = "Eigenfaces: 3-4 components, visualize via 3d scatter plot (matplotlib)/ 4d via colors, Eigenfaces: display Entwicklung von Rekonstruktionsfehler: way to find best # components, added Cluster Algo: Variationale Baysian Mixture Model after PCA"

# 2) The two TODO notes that used to live in D29 and D30 are now folded into B26 above,
#    so remove those two rows outright (not just clear their contents) while keeping every
#    other row's number the same: delete rows 29-30 (shifts 31+ up), then re-insert two
#    blank rows in their place (shifts 31+ back down) so the sheet dimension / row numbers
#    stay A1:D37.
$ws.Rows("29:30").Delete()
$ws.Rows("29:30").Insert()

# 3) Row 26 now holds much longer wrapped text, so it needs a taller row to display it.
$ws.Range("26:26").RowHeight() = 86

# 4) Reflect the edit location as the active selection, like Excel would leave it after
#    the last edited cell.
$ws.Range("B26").Select()
